# Updating units for external data (mks)
# Adds metric (m / $-per-m / kg-CO2-e-per-m-per-s) "KMS units" helper
# columns (L/M, plus a VLOOKUP demo in O/P) to the Cost and GHG sheets,
# and adds a new chart plotting the new Cost!L5:M16 series.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Cost sheet ("Cost")
# ---------------------------------------------------------------------
$cost = $wb.Worksheets.Item("Cost")

# --- Pipe cost table (rows 5:16) -> columns L (m) / M ($/m) ----------
$cost.Range("L1").Value = "KMS units"
$cost.Range("L4").Value = "m"
$cost.Range("M4").Value = "`$/m"

$cost.Range("L5").Formula = "=A5/1000"
$cost.Range("M5").Formula = "=B5"
$cost.Range("L6:L16").Formula = "=A6/1000"
$cost.Range("M6:M16").Formula = "=B6"

$cost.Range("O5").Value = 0.11
$cost.Range("P5").Formula = "=VLOOKUP(O5,L5:M16,2)"

# --- Tank cost table (rows 23:28) -> columns L (kL) / M ($) ----------
$cost.Range("L23").Formula = "=A23"
$cost.Range("M23").Formula = "=H23"
$cost.Range("L24:L28").Formula = "=A24"
$cost.Range("M24:M28").Formula = "=H24"

# --- Valve cost table (rows 33:44) -> columns L (m) / M ($) ----------
$cost.Range("L33").Formula = "=A33/1000"
$cost.Range("M33").Formula = "=B33"
$cost.Range("L34:L44").Formula = "=A34/1000"
$cost.Range("M34:M44").Formula = "=B34"

# ---------------------------------------------------------------------
# GHG sheet
# ---------------------------------------------------------------------
$ghg = $wb.Worksheets.Item("GHG")

$ghg.Range("L1").Value = "KMS units"
$ghg.Range("L4").Value = "m"
$ghg.Range("M4").Value = "kg-CO2-e/m/s"

$ghg.Range("L5").Formula = "=A5/1000"
$ghg.Range("M5").Formula = "=B5"
$ghg.Range("L6:L16").Formula = "=A6/1000"
$ghg.Range("M6:M16").Formula = "=B6"

# ---------------------------------------------------------------------
# New chart: Cost!L5:M16 (diameter in m vs $/m), smooth scatter line
# ---------------------------------------------------------------------
$chart = $cost.Shapes.AddChart2(-1, -4169).Chart
$chart.SeriesCollection.NewSeries()
$ser = $chart.SeriesCollection(1)
$ser.XValues = $cost.Range("L5:L16")
$ser.Values = $cost.Range("M5:M16")
$chart.HasLegend = $true
$chart.Legend.Position = -4152

# ---------------------------------------------------------------------
# Active sheet / selection bookkeeping (mirrors the authored workbook)
# ---------------------------------------------------------------------
$cost.Activate()
$ghg.Range("L5:M16").Select()
$cost.Range("L26").Select()
